# Auto-generated Excel COM-interop edit script
# Applies per-row leve profit recalculations across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 578.4
$ws.Range("J32").Value = 623
$ws.Range("L32").Value = 623
$ws.Range("N32").Value = -1275
$ws.Range("H64").Value = 3900.077
$ws.Range("I64").Value = 3635.7273
$ws.Range("J64").Value = 4093.9333
$ws.Range("K64").Value = 3635.7273
$ws.Range("L64").Value = 4093.9333
$ws.Range("M64").Value = -3387.7273
$ws.Range("N64").Value = -4589.933300000001
$ws.Range("H67").Value = 3900.077
$ws.Range("I67").Value = 3635.7273
$ws.Range("J67").Value = 4093.9333
$ws.Range("K67").Value = 3635.7273
$ws.Range("L67").Value = 4093.9333
$ws.Range("M67").Value = -2777.7273
$ws.Range("N67").Value = -5809.9333
$ws.Range("H95").Value = 25475
$ws.Range("J95").Value = 25475
$ws.Range("L95").Value = 25475
$ws.Range("N95").Value = -30967
$ws.Range("H123").Value = 36580
$ws.Range("J123").Value = 36580
$ws.Range("L123").Value = 36580
$ws.Range("N123").Value = -46380
$ws.Range("H135").Value = 55556588
$ws.Range("I135").Value = 22728412
$ws.Range("J135").Value = 200000580
$ws.Range("K135").Value = 204555708
$ws.Range("L135").Value = 1800005220
$ws.Range("M135").Value = -204553173
$ws.Range("N135").Value = -1800010290
$ws.Range("H137").Value = 1582.6296
$ws.Range("I137").Value = 1249.4884
$ws.Range("J137").Value = 2884.9092
$ws.Range("K137").Value = 3748.4652
$ws.Range("L137").Value = 8654.7276
$ws.Range("M137").Value = -1198.4652
$ws.Range("N137").Value = -13754.7276
$ws.Range("H138").Value = 3106.4854
$ws.Range("I138").Value = 1474.9678
$ws.Range("J138").Value = 4473.4326
$ws.Range("K138").Value = 4424.903399999999
$ws.Range("L138").Value = 13420.2978
$ws.Range("M138").Value = 715.0966000000008
$ws.Range("N138").Value = -23700.2978

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1830
$ws.Range("H32").Value = 33067.793
$ws.Range("I32").Value = 37897.773
$ws.Range("K32").Value = 37897.773
$ws.Range("M32").Value = -37610.773
$ws.Range("H52").Value = 63866.668
$ws.Range("J52").Value = 63866.668
$ws.Range("L52").Value = 63866.668
$ws.Range("N52").Value = -64502.668
$ws.Range("H110").Value = 1726.8182
$ws.Range("I110").Value = 1801.3334
$ws.Range("J110").Value = 1391.5
$ws.Range("K110").Value = 1801.3334
$ws.Range("L110").Value = 1391.5
$ws.Range("M110").Value = 243.6666
$ws.Range("N110").Value = -5481.5
$ws.Range("H115").Value = 30684
$ws.Range("J115").Value = 30684
$ws.Range("L115").Value = 30684
$ws.Range("N115").Value = -33818
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1971.7084
$ws.Range("I132").Value = 1522.193
$ws.Range("J132").Value = 3679.8667
$ws.Range("K132").Value = 4566.579
$ws.Range("L132").Value = 11039.6001
$ws.Range("M132").Value = -2036.579
$ws.Range("N132").Value = -16099.6001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1222.75
$ws.Range("I107").Value = 963.6667
$ws.Range("K107").Value = 963.6667
$ws.Range("M107").Value = 956.3333
$ws.Range("H127").Value = 55630
$ws.Range("J127").Value = 55630
$ws.Range("L127").Value = 55630
$ws.Range("N127").Value = -65550
$ws.Range("H134").Value = 46406.652
$ws.Range("I134").Value = 3128.3572
$ws.Range("J134").Value = 113728.445
$ws.Range("K134").Value = 9385.071599999999
$ws.Range("L134").Value = 341185.335
$ws.Range("M134").Value = -6850.071599999999
$ws.Range("N134").Value = -346255.335
$ws.Range("H135").Value = 38623.08
$ws.Range("J135").Value = 38623.08
$ws.Range("L135").Value = 38623.08
$ws.Range("N135").Value = -48763.08

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1277.4546
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H107").Value = 1149.3914
$ws.Range("I107").Value = 1326.3846
$ws.Range("J107").Value = 919.3
$ws.Range("K107").Value = 1326.3846
$ws.Range("L107").Value = 919.3
$ws.Range("M107").Value = 593.6153999999999
$ws.Range("N107").Value = -4759.3
$ws.Range("H132").Value = 3747.1091
$ws.Range("I132").Value = 4387.4326
$ws.Range("J132").Value = 2430.889
$ws.Range("K132").Value = 13162.2978
$ws.Range("L132").Value = 7292.667
$ws.Range("M132").Value = -10632.2978
$ws.Range("N132").Value = -12352.667
$ws.Range("H134").Value = 2889.2654
$ws.Range("I134").Value = 2377.158
$ws.Range("J134").Value = 3213.6
$ws.Range("K134").Value = 7131.474
$ws.Range("L134").Value = 9640.799999999999
$ws.Range("M134").Value = -4596.474
$ws.Range("N134").Value = -14710.8

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 24530.902
$ws.Range("J131").Value = 29948.395
$ws.Range("L131").Value = 89845.185
$ws.Range("N131").Value = -99925.185

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2745.7778
$ws.Range("I113").Value = 3800
$ws.Range("J113").Value = 2340.3076
$ws.Range("K113").Value = 3800
$ws.Range("L113").Value = 2340.3076
$ws.Range("M113").Value = -1630
$ws.Range("N113").Value = -6680.3076
$ws.Range("H132").Value = 11012.695
$ws.Range("I132").Value = 7840.2104
$ws.Range("J132").Value = 26082
$ws.Range("K132").Value = 23520.6312
$ws.Range("L132").Value = 78246
$ws.Range("M132").Value = -20990.6312
$ws.Range("N132").Value = -83306

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4549.75
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 5399.6665
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 5399.6665
$ws.Range("M12").Value = -1830
$ws.Range("N12").Value = -5739.6665
$ws.Range("H53").Value = 27400
$ws.Range("J53").Value = 27400
$ws.Range("L53").Value = 27400
$ws.Range("N53").Value = -28436
$ws.Range("H134").Value = 49000
$ws.Range("J134").Value = 49000
$ws.Range("L134").Value = 49000
$ws.Range("N134").Value = -59140
$ws.Range("H136").Value = 4057.0527
$ws.Range("I136").Value = 2283.111
$ws.Range("J136").Value = 7098.095
$ws.Range("K136").Value = 6849.333
$ws.Range("L136").Value = 21294.285
$ws.Range("M136").Value = -4299.333
$ws.Range("N136").Value = -26394.285
$ws.Range("H138").Value = 62894
$ws.Range("J138").Value = 62894
$ws.Range("L138").Value = 62894
$ws.Range("N138").Value = -73174

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3584.7144
$ws.Range("I107").Value = 2418
$ws.Range("J107").Value = 6501.5
$ws.Range("K107").Value = 7254
$ws.Range("L107").Value = 19504.5
$ws.Range("M107").Value = -5334
$ws.Range("N107").Value = -23344.5
$ws.Range("H132").Value = 2322.2727
$ws.Range("I132").Value = 1084.2307
$ws.Range("J132").Value = 4110.5557
$ws.Range("K132").Value = 3252.6921
$ws.Range("L132").Value = 12331.6671
$ws.Range("M132").Value = -722.6921000000002
$ws.Range("N132").Value = -17391.6671
$ws.Range("H136").Value = 7371.75
$ws.Range("I136").Value = 6481.5454
$ws.Range("J136").Value = 8261.954
$ws.Range("K136").Value = 19444.6362
$ws.Range("L136").Value = 24785.862
$ws.Range("M136").Value = -16894.6362
$ws.Range("N136").Value = -29885.862
